$wb = $excel.ActiveWorkbook
$wsInsert = $wb.Worksheets.Item("Insert")
$wsResult = $wb.Worksheets.Item("Result")

# --- Clean up left-over sample/test data ---------------------------------
# Row 2 and row 4 on "Insert" held sample scouting-entry strings (and the
# dependent split-out helper formulas in B:AB). Clear the whole row range
# first so no stale legacy CSE-array metadata survives, then re-enter the
# split formula as a single-cell dynamic array in column B.
$splitFormula2 = '=TRANSPOSE(TRIM(MID(SUBSTITUTE(";"&A2,";",REPT(" ",LEN(A2)+1)),ROW(INDIRECT("A1:A"&LEN(A2)-LEN(SUBSTITUTE(A2,";",""))+1))*LEN(A2)+1,LEN(A2)+1)))'
$splitFormula4 = '=TRANSPOSE(TRIM(MID(SUBSTITUTE(";"&A4,";",REPT(" ",LEN(A4)+1)),ROW(INDIRECT("A1:A"&LEN(A4)-LEN(SUBSTITUTE(A4,";",""))+1))*LEN(A4)+1,LEN(A4)+1)))'

$wsInsert.Range("A2:AB2").ClearContents()
$wsInsert.Range("B2").FormulaArray = $splitFormula2

$wsInsert.Range("A4:AB4").ClearContents()
$wsInsert.Range("B4").FormulaArray = $splitFormula4

# Stray leftover test value on "Result" sheet
$wsResult.Range("H3").ClearContents()

# --- Add the "Team Name" column ------------------------------------------
# A new column is inserted right before "Auto Start Position" on both
# sheets, pushing every later column one to the right.
$wsInsert.Columns("H:H").Insert()
$wsInsert.Range("H1").Value = "Team Name"

$wsResult.Columns("H:H").Insert()
$wsResult.Range("H1").Value = "Team Name"

# --- Selection bookmarks ---------------------------------------------------
# Update the remembered selection on each sheet, then re-activate "Insert"
# so it remains the visible/active tab (as in the original workbook).
$wsResult.Range("J26").Select()
$wsInsert.Activate()
$wsInsert.Range("A26").Select()
